$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Git command")

$data = @(
    @("git commit -m", "Save the staged changes with a commit message", 'git commit -m "message"'),
    @("git log", "Displays a log of all commits in the repository.", "git log --all"),
    @("git branch", "List all branches in the repository", "git branch"),
    @("git branch", "Creates a new branch", "git branch (branch name)"),
    @("git checkout -b", "Creates and switches to a new branch", "git checkout -b (branch name)"),
    @("git merge", "Merges another branch into the current branch", "git merge (branch name)"),
    @("git pull origin", "Fetches and merges the latest changes from a remote repository", "git pull origin (banch)")
)

$startRow = 9
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}

$ws.Range("B16").Select()
